$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland II Liga")

# Columns B (2) through AB (28) hold the actual match-record data for each
# row; column A just holds a sequential row id and must stay untouched.
$firstCol = 2   # B
$lastCol  = 28  # AB

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

# --- Simple pairwise swaps (B:AB) ---
$pairs = @(
    @(27, 28),
    @(85, 86),
    @(101, 102),
    @(108, 109),
    @(156, 157)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# --- 3-row rotation for rows 271, 272, 273 ---
# new(271) = old(273); new(272) = old(271); new(273) = old(272)
$v271 = Get-RowValues 271
$v272 = Get-RowValues 272
$v273 = Get-RowValues 273

Set-RowValues 271 $v273
Set-RowValues 272 $v271
Set-RowValues 273 $v272
